$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the "RAPID" table column to "RAPID IO" (also updates D1 header cell) ---
$ws.Range("D1").Value2 = "RAPID IO"

# --- Relabel the RAPID IO column values (col D) ---
$ws.Range("D3").Value2  = "IO-6"
$ws.Range("D4").Value2  = "IO-7"
$ws.Range("D5").Value2  = "IO-8"
$ws.Range("D29").Value2 = "SPI_CS_FLASH"
$ws.Range("D36").Value2 = "IO-4"
$ws.Range("D37").Value2 = "IO-5"
$ws.Range("D38").Value2 = "SWITCH"
$ws.Range("D43").Value2 = "USB_FS_VBUS"
$ws.Range("D44").Value2 = "HCOUT2"
$ws.Range("D45").Value2 = "USB_D-"
$ws.Range("D46").Value2 = "USB_D+"
$ws.Range("D62").Value2 = "CAN1_RX"
$ws.Range("D63").Value2 = "CAN1_TX"

# --- Font color adjustments to match new formatting ---
# D3 switches from the "gray/teal" style to the "blue" style used by D36/D37
$ws.Range("D3").Font.Color = $ws.Range("D36").Font.Color
# D38 switches from "blue" style to the plain/default style
$ws.Range("D38").Font.Color = $ws.Range("D10").Font.Color
# D29 and D44 switch from the "gray/teal" style to the plain/default style
$ws.Range("D29").Font.Color = $ws.Range("D10").Font.Color
$ws.Range("D44").Font.Color = $ws.Range("D10").Font.Color

# --- Mark G61/G64 (previously blank) and touch G62/G63 formatting ---
$ws.Range("G61").WrapText = $true
$ws.Range("G64").WrapText = $true
$ws.Range("G62").WrapText = $true
$ws.Range("G63").WrapText = $true

# --- Sheet view: zoom back to 100%, scroll to top-left, select F5 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("F5").Select()
